$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Remove the entire row 36 (Caso -461, Independencia 2796, San Telmo).
# All subsequent rows shift up by one, matching the target diff.
$ws.Rows.Item(36).Delete()
